$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 2 corresponds to the 3b5cfbe7... file.
# Update "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
$wsZhCn.Range("E2").Value = "2016-03-24 08:26:57"
$wsZhCn.Range("H2").Value = "2016-03-24 08:27:22"

# de-de sheet: row 2 corresponds to the 3b5cfbe7... file.
$wsDeDe.Range("E2").Value = "2016-03-24 08:27:02"
$wsDeDe.Range("H2").Value = "2016-03-24 08:27:29"
